$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1586

$ws.Range("H19").Value = 975.8
$ws.Range("I19").Value = 1549.6666
$ws.Range("J19").Value = 115
$ws.Range("K19").Value = 1549.6666
$ws.Range("L19").Value = 115
$ws.Range("M19").Value = -1374.6666
$ws.Range("N19").Value = -465

$ws.Range("H55").Value = 430.66666
$ws.Range("I55").Value = 310
$ws.Range("J55").Value = 672
$ws.Range("K55").Value = 310
$ws.Range("L55").Value = 672
$ws.Range("M55").Value = -96
$ws.Range("N55").Value = -1100

$ws.Range("H112").Value = 2540.3462
$ws.Range("J112").Value = 2574
$ws.Range("L112").Value = 7722
$ws.Range("N112").Value = -9938

$ws.Range("H141").Value = 3550.125
$ws.Range("I141").Value = 3952.8572
$ws.Range("K141").Value = 11858.5716
$ws.Range("M141").Value = -6678.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 983.2
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 5000
$ws.Range("M57").Value = -4516

$ws.Range("H74").Value = 1744.8572
$ws.Range("I74").Value = 1785.6666
$ws.Range("K74").Value = 1785.6666
$ws.Range("M74").Value = -911.6666

$ws.Range("H77").Value = 1744.8572
$ws.Range("I77").Value = 1785.6666
$ws.Range("K77").Value = 8928.333000000001
$ws.Range("M77").Value = -4560.333000000001

$ws.Range("H97").Value = 348.5
$ws.Range("I97").Value = 348.5
$ws.Range("K97").Value = 348.5
$ws.Range("M97").Value = 147.5

$ws.Range("H110").Value = 1536.3334
$ws.Range("I110").Value = 546.8570999999999
$ws.Range("K110").Value = 546.8570999999999
$ws.Range("M110").Value = 1498.1429

$ws.Range("H122").Value = 2401.5715
$ws.Range("I122").Value = 2339.1052
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 7017.3156
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -4567.3156
$ws.Range("N122").Value = -13885

$ws.Range("H132").Value = 1411.0646
$ws.Range("I132").Value = 1411.0646
$ws.Range("K132").Value = 4233.1938
$ws.Range("M132").Value = -1703.1938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1299.7693
$ws.Range("I20").Value = 573.6667
$ws.Range("K20").Value = 573.6667
$ws.Range("M20").Value = -326.6667

$ws.Range("H99").Value = 2120
$ws.Range("I99").Value = 2050
$ws.Range("K99").Value = 2050
$ws.Range("M99").Value = -552

$ws.Range("H105").Value = 2043.375
$ws.Range("I105").Value = 2057.8333
$ws.Range("K105").Value = 2057.8333
$ws.Range("M105").Value = -310.8332999999998

$ws.Range("H134").Value = 7099.8667
$ws.Range("I134").Value = 7045.273
$ws.Range("K134").Value = 21135.819
$ws.Range("M134").Value = -18600.819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 88.882355
$ws.Range("I7").Value = 104.63636
$ws.Range("K7").Value = 104.63636
$ws.Range("M7").Value = 8.363640000000004

$ws.Range("H31").Value = 2517
$ws.Range("I31").Value = 2210.6667
$ws.Range("J31").Value = 2700.8
$ws.Range("K31").Value = 2210.6667
$ws.Range("L31").Value = 2700.8
$ws.Range("M31").Value = -1915.6667
$ws.Range("N31").Value = -3290.8

$ws.Range("H34").Value = 2517
$ws.Range("I34").Value = 2210.6667
$ws.Range("J34").Value = 2700.8
$ws.Range("K34").Value = 2210.6667
$ws.Range("L34").Value = 2700.8
$ws.Range("M34").Value = -2008.6667
$ws.Range("N34").Value = -3104.8

$ws.Range("H58").Value = 2298
$ws.Range("I58").Value = 2252.3845
$ws.Range("J58").Value = 2495.6667
$ws.Range("K58").Value = 2252.3845
$ws.Range("L58").Value = 2495.6667
$ws.Range("M58").Value = -2049.3845
$ws.Range("N58").Value = -2901.6667

$ws.Range("H136").Value = 2298
$ws.Range("I136").Value = 2252.3845
$ws.Range("J136").Value = 2495.6667
$ws.Range("K136").Value = 6757.1535
$ws.Range("L136").Value = 7487.000100000001
$ws.Range("M136").Value = -4207.1535
$ws.Range("N136").Value = -12587.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3954.9092
$ws.Range("I140").Value = 3722.7778
$ws.Range("K140").Value = 11168.3334
$ws.Range("M140").Value = -5988.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 524.5
$ws.Range("J2").Value = 550
$ws.Range("L2").Value = 550
$ws.Range("N2").Value = -776

$ws.Range("H80").Value = 3590.3103
$ws.Range("I80").Value = 2126.4
$ws.Range("J80").Value = 3895.2917
$ws.Range("K80").Value = 2126.4
$ws.Range("L80").Value = 3895.2917
$ws.Range("M80").Value = -1128.4
$ws.Range("N80").Value = -5891.2917

$ws.Range("H83").Value = 3590.3103
$ws.Range("I83").Value = 2126.4
$ws.Range("J83").Value = 3895.2917
$ws.Range("K83").Value = 10632
$ws.Range("L83").Value = 19476.4585
$ws.Range("M83").Value = -5640
$ws.Range("N83").Value = -29460.4585

$ws.Range("H122").Value = 8933294
$ws.Range("I122").Value = 8933294
$ws.Range("K122").Value = 26799882
$ws.Range("M122").Value = -26797432

$ws.Range("H139").Value = 116666.664
$ws.Range("J139").Value = 116666.664
$ws.Range("L139").Value = 116666.664
$ws.Range("N139").Value = -126946.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4749.25
$ws.Range("I7").Value = 3004
$ws.Range("J7").Value = 4998.5713
$ws.Range("K7").Value = 3004
$ws.Range("L7").Value = 4998.5713
$ws.Range("M7").Value = -2892
$ws.Range("N7").Value = -5222.5713

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H126").Value = 4749.25
$ws.Range("I126").Value = 3004
$ws.Range("J126").Value = 4998.5713
$ws.Range("K126").Value = 9012
$ws.Range("L126").Value = 14995.7139
$ws.Range("M126").Value = -6542
$ws.Range("N126").Value = -19935.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 208.75
$ws.Range("I81").Value = 208.75
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 417.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 643.5
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 208.75
$ws.Range("I84").Value = 208.75
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 2087.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 3216.5
$ws.Range("N84").ClearContents()
